$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1) Insert a new "2024" year-header row before what is currently the Jan-2024 row (116) ---
$ws.Rows(116).Insert()

# New row inherits formatting from the row above (matches the normal data-row look)...
$ws.Range("A115:K115").Copy()
$ws.Range("A116:K116").PasteSpecial(-4122)

# ...except column A, which should look like the other year-header cells (e.g. A103 = "2023")
$ws.Range("A103").Copy()
$ws.Range("A116").PasteSpecial(-4122)
$ws.Range("A116").Value = "'2024"

# Restore the calculated-column formula in the new row (format-only paste clears formulas)
$ws.Range("G116").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# --- 2) Grow Table1 so it keeps covering the whole data range (A8:K138 -> A8:K139) ---
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A8:K139"))

# Restore the calculated-column formula on the table's trailing row too
$ws.Range("G139").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# --- 3) Fill in leave data that was already present for existing rows (Aug-Dec 2023) ---
# August 2023 (row 111) - SL(1-0-0) already marked; add EARNED 1.25
$ws.Range("C111").Value = 1.25

# September 2023 (row 112) - add EARNED 1.25
$ws.Range("C112").Value = 1.25

# October 2023 (row 113) - mark SL(1-0-0), EARNED 1.25, balance-forward flag, and a remark date
$ws.Range("B113").Value = "SL(1-0-0)"
$ws.Range("C113").Value = 1.25
$ws.Range("H113").Value = 1
$ws.Range("K111").Copy()
$ws.Range("K113").PasteSpecial(-4122)
$ws.Range("K113").Value = 45206

# November 2023 (row 114) - add EARNED 1.25
$ws.Range("C114").Value = 1.25

# December 2023 (row 115) - mark SL(1-0-0), balance-forward flag, and a remark date
$ws.Range("B115").Value = "SL(1-0-0)"
$ws.Range("H115").Value = 1
$ws.Range("K111").Copy()
$ws.Range("K115").PasteSpecial(-4122)
$ws.Range("K115").Value = 45261

# --- 4) Reflect where the user's selection/scroll ended up when the file was saved ---
$ws.Application.ActiveWindow.ScrollRow = 101
$ws.Range("K115").Select()
